$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the "actual run" label (shared string) to "actual batch" and add
#    a new batch-size row header ("batch") plus a "no beef or lamb" section
#    title string used further down.
# ---------------------------------------------------------------------------
$ws.Range("A163").Value = "actual batch"

# ---------------------------------------------------------------------------
# 2) Row 163 ("actual batch"): fill in batch index numbers 1..9 across B:J
# ---------------------------------------------------------------------------
$batchNumbers = @(1,2,3,4,5,6,7,8,9)
for ($i = 0; $i -lt $batchNumbers.Length; $i++) {
    $ws.Cells.Item(163, 2 + $i).Value = $batchNumbers[$i]
}

# ---------------------------------------------------------------------------
# 3) Row 164 ("test runs"): new per-batch sample-size values
# ---------------------------------------------------------------------------
$testRuns = @(3011,3149,3054,3093,3154,3077,3026,3012,2000)
for ($i = 0; $i -lt $testRuns.Length; $i++) {
    $ws.Cells.Item(164, 2 + $i).Value = $testRuns[$i]
}

# ---------------------------------------------------------------------------
# 4) Row 165 ("cost £"): new per-batch cost values
# ---------------------------------------------------------------------------
$costVals = @(27241,30681,31253,29389,30123,30477,29196,28656,20250)
for ($i = 0; $i -lt $costVals.Length; $i++) {
    $ws.Cells.Item(165, 2 + $i).Value = $costVals[$i]
}

# ---------------------------------------------------------------------------
# 5) Row 166 ("emissions kg"): new per-batch emissions values
# ---------------------------------------------------------------------------
$emissionsVals = @(24796,27340,27585,28174,28021,27233,25906,24330,17068)
for ($i = 0; $i -lt $emissionsVals.Length; $i++) {
    $ws.Cells.Item(166, 2 + $i).Value = $emissionsVals[$i]
}

# ---------------------------------------------------------------------------
# 6) Row 167 ("food waste cals"): new per-batch food-waste values
# ---------------------------------------------------------------------------
$foodWasteVals = @(1128873,826028,747515,1252003,1043889,852864,768373,1098457,607663)
for ($i = 0; $i -lt $foodWasteVals.Length; $i++) {
    $ws.Cells.Item(167, 2 + $i).Value = $foodWasteVals[$i]
}

# ---------------------------------------------------------------------------
# 7) Row 168 ("cost/(days x people)"): cost / test-runs formulas per batch
#    (column B keeps its historical parenthesised style)
# ---------------------------------------------------------------------------
$ws.Range("B168").Formula = "=(27241/3011)"
$ws.Range("C168").Formula = "=30681/3149"
$ws.Range("D168").Formula = "=31253/3054"
$ws.Range("E168").Formula = "=29389/3093"
$ws.Range("F168").Formula = "=30123/3154"
$ws.Range("G168").Formula = "=30477/3077"
$ws.Range("H168").Formula = "=29196/3026"
$ws.Range("I168").Formula = "=28656/3012"
$ws.Range("J168").Formula = "=20250/2000"

# ---------------------------------------------------------------------------
# 8) Row 169 ("emissions /(days x people)")
# ---------------------------------------------------------------------------
$ws.Range("B169").Formula = "=(24796/3011)"
$ws.Range("C169").Formula = "=27340/3149"
$ws.Range("D169").Formula = "=27585/3054"
$ws.Range("E169").Formula = "=28174/3093"
$ws.Range("F169").Formula = "=28021/3154"
$ws.Range("G169").Formula = "=27233/3077"
$ws.Range("H169").Formula = "=25906/3026"
$ws.Range("I169").Formula = "=24330/3012"
$ws.Range("J169").Formula = "=17068/2000"

# ---------------------------------------------------------------------------
# 9) Row 170 ("food waste /(days x people)")
# ---------------------------------------------------------------------------
$ws.Range("B170").Formula = "=(1128873/3011)"
$ws.Range("C170").Formula = "=826028/3149"
$ws.Range("D170").Formula = "=747515/3054"
$ws.Range("E170").Formula = "=1252003/3093"
$ws.Range("F170").Formula = "=1043889/3154"
$ws.Range("G170").Formula = "=852864/3077"
$ws.Range("H170").Formula = "=768373/3026"
$ws.Range("I170").Formula = "=1098457/3012"
$ws.Range("J170").Formula = "=607663/2000"

# ---------------------------------------------------------------------------
# 10) New section starting at row 172: "no beef or lamb" batch comparing the
#     first batch (B) against the last one (J) only.
# ---------------------------------------------------------------------------
$ws.Range("A172").Value = "no beef or lamb"

$ws.Range("A173").Value = "batch"
$ws.Range("B173").Value = 1
$ws.Range("J173").Value = 9

$ws.Range("A174").Value = "Matrix size"
$ws.Range("B174").Value = 3011
$ws.Range("J174").Value = 2000

$ws.Range("A175").Value = "cost £"
$ws.Range("B175").Value = 27665
$ws.Range("J175").Value = 20337

$ws.Range("A176").Value = "emissions kg"
$ws.Range("B176").Value = 21977
$ws.Range("J176").Value = 15733

$ws.Range("A177").Value = "food waste sum nutrients"
$ws.Range("B177").Value = 830515
$ws.Range("J177").Value = 600977

$ws.Range("A178").Value = "cost/(days x people)"
$ws.Range("B178").Formula = "=27665/3011"
$ws.Range("J178").Formula = "=20337/2000"

$ws.Range("A179").Value = "emissions /(days x people)"
$ws.Range("B179").Formula = "=21977/3011"
$ws.Range("J179").Formula = "=15733/2000"

$ws.Range("A180").Value = "food waste /(days x people)"
$ws.Range("B180").Formula = "=830515/3011"
$ws.Range("J180").Formula = "=600977/2000"

# ---------------------------------------------------------------------------
# 11) Column A is now much wider (holds the longer new row labels) - resize
#     it to match, dropping the old best-fit auto width.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 22.5

# ---------------------------------------------------------------------------
# 12) Update the view: scrolled further down, new active selection.
# ---------------------------------------------------------------------------
$ws.Range("I178").Select()
